$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "version" (lesson name/version) column in column I.
# Header in row 2, matching the existing bold header row style.
$ws.Range("I2").Value = "version"
$ws.Range("I2").Font.Bold = $true
$ws.Range("I2").Font.Name = "Calibri"
$ws.Range("I2").Font.Size = 11

# Fill I3:I18 with a repeating A / B / C pattern (one per data row).
$values = @("A","B","C")
for ($row = 3; $row -le 18; $row++) {
    $idx = ($row - 3) % 3
    $cell = $ws.Range("I" + $row)
    $cell.Value = $values[$idx]
    $cell.Font.Bold = $false
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
}

# Page setup: portrait, paper size 9 (A4), matching the author's printer settings.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Match the final cursor/selection position left by the author.
$null = $ws.Range("H20").Select()
